$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price/volume columns remain plain text so values like "1.001" or
# "30.588.64" are not reinterpreted as numbers/dates by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "30.588.64"
$ws.Range("D3").Value = "1.922.94"
$ws.Range("E3").Value = "  -0.56%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "247.26"
$ws.Range("E5").Value = "  +2.76%  "
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("D7").Value = "0.4730"
$ws.Range("E7").Value = "  -0.49%  "
$ws.Range("D8").Value = "0.2913"
$ws.Range("E8").Value = "  +1.20%  "
$ws.Range("D9").Value = "0.06822"
$ws.Range("E9").Value = "  +2.62%  "
$ws.Range("D10").Value = "104.94"
$ws.Range("E10").Value = "  -2.34%  "
$ws.Range("D11").Value = "18.37"
$ws.Range("E11").Value = "  -4.11%  "
$ws.Range("D12").Value = "1.929.26"
$ws.Range("E12").Value = "  -0.19%  "
$ws.Range("D13").Value = "0.07721"
$ws.Range("E13").Value = "  +1.35%  "
$ws.Range("D14").Value = "5.322"
$ws.Range("E14").Value = "  +3.03%  "
$ws.Range("D15").Value = "0.6697"
$ws.Range("E15").Value = "  +0.90%  "
$ws.Range("D16").Value = "290.37"
$ws.Range("E16").Value = "  -5.60%  "
$ws.Range("D17").Value = "30.609.41"
$ws.Range("E17").Value = "  +0.35%  "
$ws.Range("D18").Value = "0.000007627"
$ws.Range("E18").Value = "  +0.34%  "
$ws.Range("D19").Value = "1.001"
$ws.Range("E19").Value = "  +0.10%  "
$ws.Range("D20").Value = "12.94"
$ws.Range("E20").Value = "  -0.62%  "
$ws.Range("D21").Value = "5.551"
$ws.Range("E21").Value = "  +4.72%  "
$ws.Range("D22").Value = "2.172.36"
$ws.Range("E22").Value = "  -0.26%  "
$ws.Range("D23").Value = "1.002"
$ws.Range("E23").Value = "  +0.15%  "
$ws.Range("D24").Value = "6.469"
$ws.Range("E24").Value = "  +2.64%  "
$ws.Range("D25").Value = "9.511"
$ws.Range("D26").Value = "167.21"
$ws.Range("E26").Value = "  -0.56%  "
$ws.Range("D27").Value = "20.77"
$ws.Range("E27").Value = "  +2.40%  "
$ws.Range("D28").Value = "2.136"
$ws.Range("E28").Value = "  +4.11%  "
$ws.Range("D29").Value = "0.1069"
$ws.Range("E29").Value = "  -3.59%  "
$ws.Range("D30").Value = "1.407"
$ws.Range("E30").Value = "  +2.70%  "
$ws.Range("D31").Value = "4.192"
$ws.Range("E31").Value = "  +2.06%  "
$ws.Range("D32").Value = "4.050"
$ws.Range("E32").Value = "  +2.81%  "
$ws.Range("D33").Value = "0.05023"
$ws.Range("E33").Value = "  -0.13%  "
$ws.Range("D34").Value = "0.7335"
$ws.Range("E34").Value = "  -1.38%  "
$ws.Range("D35").Value = "1.144"
$ws.Range("E35").Value = "  -1.11%  "
$ws.Range("D36").Value = "0.02058"
$ws.Range("E36").Value = "  +4.67%  "
$ws.Range("D37").Value = "0.9996"
$ws.Range("E37").Value = "  +0.09%  "
$ws.Range("D38").Value = "2.738"
$ws.Range("E38").Value = "  -0.66%  "
$ws.Range("E39").Value = "  -0.37%  "
$ws.Range("D40").Value = "111.82"
$ws.Range("E40").Value = "  +3.75%  "
$ws.Range("D41").Value = "2.041"
$ws.Range("E41").Value = "  -0.41%  "
$ws.Range("D42").Value = "0.4438"
$ws.Range("E42").Value = "  +5.67%  "
$ws.Range("D43").Value = "0.8717"
$ws.Range("E43").Value = "  -1.12%  "
$ws.Range("D44").Value = "5.885"
$ws.Range("E44").Value = "  +1.12%  "
$ws.Range("E45").Value = "  +0.10%  "
$ws.Range("D46").Value = "67.76"
$ws.Range("E46").Value = "  -3.84%  "
$ws.Range("D47").Value = "7.274"
$ws.Range("E47").Value = "  -0.26%  "
$ws.Range("D48").Value = "9.401"
$ws.Range("E48").Value = "  +1.23%  "
$ws.Range("D49").Value = "0.1250"
$ws.Range("E49").Value = "  +2.82%  "
$ws.Range("D50").Value = "47.89"
$ws.Range("E50").Value = "  +12.94%  "
$ws.Range("D51").Value = "35.18"
$ws.Range("E51").Value = "  +0.75%  "
